$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new metadata row (row 2) ---
# Columns: A=identifier, C=title, D=(blank), E=levelOfDescription,
#          F=extentAndMedium, G=notes, H=(blank)
$ws.Range("A2").Value = "MCH133"
$ws.Range("C2").Value = "VISTA U NIVERITY, THE HISTOORY OF THOMAS MTOBI, UWC CAMPUS MAIL"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: Cabinet 1B | GRAP COUNT NUMER: NONE"

# D2 and H2 stay empty, but still carry the row's styling (set below).

# --- Apply the row's cell formatting (10pt Calibri, theme text color) ---
$dataCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $dataCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}

# --- Select row2 / freeze header row, matching the updated sheet view ---
$ws.Range("A2:I2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
